$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.043.38"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.354.04"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.27"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.54"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  +3.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.64"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.03"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "2.774.26"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").Value = "57.928.18"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "2.359.65"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "330.40"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.76"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.26"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("E27").Value = "  -5.38%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.74"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "0.0₃0737"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.42"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.434"
$ws.Range("E38").Value = "  +15.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.29"
$ws.Range("E39").Value = "  +3.16%  "
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.67"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.48"
$ws.Range("E42").Value = "  -4.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "287.66"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("B44").Value = "Polygon"
$ws.Range("C44").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.426"
$ws.Range("E44").Value = "  +10.85%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0955"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0515"
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.566"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.56"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.08"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  +0.11%  "
